$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Stok_awal (D) and Harga (G) values for rows 2-51 ---
$ws.Range("D2").Value = 125
$ws.Range("G2").Value = 3000
$ws.Range("G3").Value = 10000
$ws.Range("G4").Value = 7000
$ws.Range("D5").Value = 125
$ws.Range("G5").Value = 15000
$ws.Range("G6").Value = 10000
$ws.Range("D7").Value = 150
$ws.Range("G7").Value = 8000
$ws.Range("G8").Value = 6000
$ws.Range("G9").Value = 5000
$ws.Range("G10").Value = 15000
$ws.Range("G11").Value = 15000
$ws.Range("G12").Value = 5000
$ws.Range("G13").Value = 20000
$ws.Range("G14").Value = 15000
$ws.Range("G15").Value = 10000
$ws.Range("G16").Value = 8000
$ws.Range("G17").Value = 20000
$ws.Range("G18").Value = 15000
$ws.Range("G19").Value = 10000
$ws.Range("G20").Value = 15000
$ws.Range("G21").Value = 30000
$ws.Range("G22").Value = 15000
$ws.Range("G23").Value = 20000
$ws.Range("G24").Value = 25000
$ws.Range("G25").Value = 10000
$ws.Range("G26").Value = 10000
$ws.Range("G27").Value = 20000
$ws.Range("G28").Value = 8000
$ws.Range("G29").Value = 15000
$ws.Range("G30").Value = 10000
$ws.Range("G31").Value = 15000
$ws.Range("G32").Value = 30000
$ws.Range("G33").Value = 10000
$ws.Range("G34").Value = 30000
$ws.Range("G35").Value = 40000
$ws.Range("G36").Value = 50000
$ws.Range("G37").Value = 15000
$ws.Range("G38").Value = 30000
$ws.Range("G39").Value = 30000
$ws.Range("G40").Value = 40000
$ws.Range("G41").Value = 30000
$ws.Range("G42").Value = 50000
$ws.Range("G43").Value = 30000
$ws.Range("G44").Value = 40000
$ws.Range("G45").Value = 20000
$ws.Range("G46").Value = 30000
$ws.Range("G47").Value = 40000
$ws.Range("G48").Value = 25000
$ws.Range("G49").Value = 40000
$ws.Range("G50").Value = 30000
$ws.Range("G51").Value = 50000
# --- Remove the now-obsolete extra rows 52-58 ---
$ws.Range("A52:G58").EntireRow.Delete()

# --- Column widths (best-fit sizing for the now-visible data columns) ---
$ws.Columns("B").ColumnWidth = 11.666666666666666
$ws.Columns("C").ColumnWidth = 18.5
$ws.Columns("D").ColumnWidth = 9.166666666666666
$ws.Columns("E").ColumnWidth = 11.166666666666666

# --- Restore the selection / scroll position used when the file was saved ---
$ws.Range("I39").Select()
